$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A2").Value = "PRE_TEK49_RES_1950"
$ws.Range("A3").Value = "PRE_TEK49_RES_1940"
$ws.Range("A7").Value = "TEK69_RES_1976"
$ws.Range("A8").Value = "TEK69_RES_1986"

$ws.Columns.Item(1).ColumnWidth = 26.15

$ws.Range("D10").Select()
